$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 74: change condition (C74) and name (F74) to the new CIMSOURCE_OPTIMALIZECONFIG case
$ws.Range('C74').Value = 'opt_cfg_emi_frac=''0.5'''
$ws.Range('F74').Value = 'CIMSOURCE_OPTIMALIZECONFIG'

# Add new test-case rows 77-103 for condition with different type combinations
$ws.Range('A77').Value = 'iems-connector-test-mysql-Float_varchar-LocalDateTime_datetime-1'
$ws.Range('B77').Value = 'good request, data retrieved (no schema check)'
$ws.Range('C77').Value = 'motor_current_percent=''92.0'' and updateTime<''2021-05-21T14:22:43'''
$ws.Range('F77').Value = 'HeatPumpKpiData'
$ws.Range('H77').Value = 0
$ws.Range('I77').Value = 0
$ws.Range('L77').Value = 0

$ws.Range('A78').Value = 'iems-connector-test-mysql-Float_varchar-String_varchar-1'
$ws.Range('B78').Value = 'good request, data retrieved (no schema check)'
$ws.Range('C78').Value = 'motor_current_percent=''92.0'' and deviceName=''1#制冷机'''
$ws.Range('F78').Value = 'HeatPumpKpiData'
$ws.Range('H78').Value = 0
$ws.Range('I78').Value = 0
$ws.Range('L78').Value = 0

$ws.Range('A79').Value = 'iems-connector-test-mysql-Float_varchar-Long_int-1'
$ws.Range('B79').Value = 'good request, data retrieved (no schema check)'
$ws.Range('C79').Value = 'motor_current_percent=''92.0'' and id=''38'''
$ws.Range('F79').Value = 'HeatPumpKpiData'
$ws.Range('H79').Value = 0
$ws.Range('I79').Value = 0
$ws.Range('L79').Value = 0

$ws.Range('A80').Value = 'iems-connector-test-mysql-LocalDateTime_datetime-String_varchar-1'
$ws.Range('B80').Value = 'good request, data retrieved (no schema check)'
$ws.Range('C80').Value = 'updateTime<''2021-05-21T14:22:43'' and deviceName=''1#制冷机'''
$ws.Range('F80').Value = 'HeatPumpKpiData'
$ws.Range('H80').Value = 0
$ws.Range('I80').Value = 0
$ws.Range('L80').Value = 0

$ws.Range('A81').Value = 'iems-connector-test-mysql-LocalDateTime_datetime-Long_int-1'
$ws.Range('B81').Value = 'good request, data retrieved (no schema check)'
$ws.Range('C81').Value = 'updateTime<''2021-05-21T14:22:43'' and id=''38'''
$ws.Range('F81').Value = 'HeatPumpKpiData'
$ws.Range('H81').Value = 0
$ws.Range('I81').Value = 0
$ws.Range('L81').Value = 0

$ws.Range('A82').Value = 'iems-connector-test-mysql-String_varchar-Long_int-1'
$ws.Range('B82').Value = 'good request, data retrieved (no schema check)'
$ws.Range('C82').Value = 'deviceName=''1#制冷机'' and id=''38'''
$ws.Range('F82').Value = 'HeatPumpKpiData'
$ws.Range('H82').Value = 0
$ws.Range('I82').Value = 0
$ws.Range('L82').Value = 0

$ws.Range('A83').Value = 'iems-connector-test-mysql-Float_varchar-LocalDateTime_datetime-String_varchar-1'
$ws.Range('B83').Value = 'good request, data retrieved (no schema check)'
$ws.Range('C83').Value = 'motor_current_percent=''92.0'' and updateTime<''2021-05-21T14:22:43'' and deviceName=''1#制冷机'''
$ws.Range('F83').Value = 'HeatPumpKpiData'
$ws.Range('H83').Value = 0
$ws.Range('I83').Value = 0
$ws.Range('L83').Value = 0

$ws.Range('A84').Value = 'iems-connector-test-mysql-Float_varchar-LocalDateTime_datetime-Long_int-1'
$ws.Range('B84').Value = 'good request, data retrieved (no schema check)'
$ws.Range('C84').Value = 'motor_current_percent=''92.0'' and updateTime<''2021-05-21T14:22:43'' and id=''38'''
$ws.Range('F84').Value = 'HeatPumpKpiData'
$ws.Range('H84').Value = 0
$ws.Range('I84').Value = 0
$ws.Range('L84').Value = 0

$ws.Range('A85').Value = 'iems-connector-test-mysql-Float_varchar-String_varchar-Long_int-1'
$ws.Range('B85').Value = 'good request, data retrieved (no schema check)'
$ws.Range('C85').Value = 'motor_current_percent=''92.0'' and deviceName=''1#制冷机'' and id=''38'''
$ws.Range('F85').Value = 'HeatPumpKpiData'
$ws.Range('H85').Value = 0
$ws.Range('I85').Value = 0
$ws.Range('L85').Value = 0

$ws.Range('A86').Value = 'iems-connector-test-mysql-LocalDateTime_datetime-String_varchar-Long_int-1'
$ws.Range('B86').Value = 'good request, data retrieved (no schema check)'
$ws.Range('C86').Value = 'updateTime<''2021-05-21T14:22:43'' and deviceName=''1#制冷机'' and id=''38'''
$ws.Range('F86').Value = 'HeatPumpKpiData'
$ws.Range('H86').Value = 0
$ws.Range('I86').Value = 0
$ws.Range('L86').Value = 0

$ws.Range('A87').Value = 'iems-connector-test-mysql-Float_varchar-LocalDateTime_datetime-String_varchar-Long_int-1'
$ws.Range('B87').Value = 'good request, data retrieved (no schema check)'
$ws.Range('C87').Value = 'motor_current_percent=''92.0'' and updateTime<''2021-05-21T14:22:43'' and deviceName=''1#制冷机'' and id=''38'''
$ws.Range('F87').Value = 'HeatPumpKpiData'
$ws.Range('H87').Value = 0
$ws.Range('I87').Value = 0
$ws.Range('L87').Value = 0

$ws.Range('A88').Value = 'iems-connector-test-mysql-Float_varchar-LocalDateTime_datetime-2'
$ws.Range('B88').Value = 'good request, data retrieved (no schema check)'
$ws.Range('C88').Value = 'motor_current_percent=''92.0'' or updateTime<''2021-05-21T14:22:43'''
$ws.Range('F88').Value = 'HeatPumpKpiData'
$ws.Range('H88').Value = 0
$ws.Range('I88').Value = 0
$ws.Range('L88').Value = 0

$ws.Range('A89').Value = 'iems-connector-test-mysql-Float_varchar-String_varchar-2'
$ws.Range('B89').Value = 'good request, data retrieved (no schema check)'
$ws.Range('C89').Value = 'motor_current_percent=''92.0'' or deviceName=''1#制冷机'''
$ws.Range('F89').Value = 'HeatPumpKpiData'
$ws.Range('H89').Value = 0
$ws.Range('I89').Value = 0
$ws.Range('L89').Value = 0

$ws.Range('A90').Value = 'iems-connector-test-mysql-Float_varchar-Long_int-2'
$ws.Range('B90').Value = 'good request, data retrieved (no schema check)'
$ws.Range('C90').Value = 'motor_current_percent=''92.0'' or id=''38'''
$ws.Range('F90').Value = 'HeatPumpKpiData'
$ws.Range('H90').Value = 0
$ws.Range('I90').Value = 0
$ws.Range('L90').Value = 0

$ws.Range('A91').Value = 'iems-connector-test-mysql-LocalDateTime_datetime-String_varchar-2'
$ws.Range('B91').Value = 'good request, data retrieved (no schema check)'
$ws.Range('C91').Value = 'updateTime<''2021-05-21T14:22:43'' or deviceName=''1#制冷机'''
$ws.Range('F91').Value = 'HeatPumpKpiData'
$ws.Range('H91').Value = 0
$ws.Range('I91').Value = 0
$ws.Range('L91').Value = 0

$ws.Range('A92').Value = 'iems-connector-test-mysql-LocalDateTime_datetime-Long_int-2'
$ws.Range('B92').Value = 'good request, data retrieved (no schema check)'
$ws.Range('C92').Value = 'updateTime<''2021-05-21T14:22:43'' or id=''38'''
$ws.Range('F92').Value = 'HeatPumpKpiData'
$ws.Range('H92').Value = 0
$ws.Range('I92').Value = 0
$ws.Range('L92').Value = 0

$ws.Range('A93').Value = 'iems-connector-test-mysql-String_varchar-Long_int-2'
$ws.Range('B93').Value = 'good request, data retrieved (no schema check)'
$ws.Range('C93').Value = 'deviceName=''1#制冷机'' or id=''38'''
$ws.Range('F93').Value = 'HeatPumpKpiData'
$ws.Range('H93').Value = 0
$ws.Range('I93').Value = 0
$ws.Range('L93').Value = 0

$ws.Range('A94').Value = 'iems-connector-test-mysql-Float_varchar-LocalDateTime_datetime-String_varchar-2'
$ws.Range('B94').Value = 'good request, data retrieved (no schema check)'
$ws.Range('C94').Value = 'motor_current_percent=''92.0'' or updateTime<''2021-05-21T14:22:43'' or deviceName=''1#制冷机'''
$ws.Range('F94').Value = 'HeatPumpKpiData'
$ws.Range('H94').Value = 0
$ws.Range('I94').Value = 0
$ws.Range('L94').Value = 0

$ws.Range('A95').Value = 'iems-connector-test-mysql-Float_varchar-LocalDateTime_datetime-Long_int-2'
$ws.Range('B95').Value = 'good request, data retrieved (no schema check)'
$ws.Range('C95').Value = 'motor_current_percent=''92.0'' or updateTime<''2021-05-21T14:22:43'' or id=''38'''
$ws.Range('F95').Value = 'HeatPumpKpiData'
$ws.Range('H95').Value = 0
$ws.Range('I95').Value = 0
$ws.Range('L95').Value = 0

$ws.Range('A96').Value = 'iems-connector-test-mysql-Float_varchar-String_varchar-Long_int-2'
$ws.Range('B96').Value = 'good request, data retrieved (no schema check)'
$ws.Range('C96').Value = 'motor_current_percent=''92.0'' or deviceName=''1#制冷机'' or id=''38'''
$ws.Range('F96').Value = 'HeatPumpKpiData'
$ws.Range('H96').Value = 0
$ws.Range('I96').Value = 0
$ws.Range('L96').Value = 0

$ws.Range('A97').Value = 'iems-connector-test-mysql-LocalDateTime_datetime-String_varchar-Long_int-2'
$ws.Range('B97').Value = 'good request, data retrieved (no schema check)'
$ws.Range('C97').Value = 'updateTime<''2021-05-21T14:22:43'' or deviceName=''1#制冷机'' or id=''38'''
$ws.Range('F97').Value = 'HeatPumpKpiData'
$ws.Range('H97').Value = 0
$ws.Range('I97').Value = 0
$ws.Range('L97').Value = 0

$ws.Range('A98').Value = 'iems-connector-test-mysql-Float_varchar-LocalDateTime_datetime-String_varchar-Long_int-2'
$ws.Range('B98').Value = 'good request, data retrieved (no schema check)'
$ws.Range('C98').Value = 'motor_current_percent=''92.0'' or updateTime<''2021-05-21T14:22:43'' or deviceName=''1#制冷机'' or id=''38'''
$ws.Range('F98').Value = 'HeatPumpKpiData'
$ws.Range('H98').Value = 0
$ws.Range('I98').Value = 0
$ws.Range('L98').Value = 0

$ws.Range('A99').Value = 'iems-connector-test-mysql-Float_varchar-LocalDateTime_datetime-String_varchar-3'
$ws.Range('B99').Value = 'good request, data retrieved (no schema check)'
$ws.Range('C99').Value = 'motor_current_percent=''92.0'' and updateTime<''2021-05-21T14:22:43'' or deviceName=''1#制冷机'''
$ws.Range('F99').Value = 'HeatPumpKpiData'
$ws.Range('H99').Value = 0
$ws.Range('I99').Value = 0
$ws.Range('L99').Value = 0

$ws.Range('A100').Value = 'iems-connector-test-mysql-Float_varchar-LocalDateTime_datetime-Long_int-3'
$ws.Range('B100').Value = 'good request, data retrieved (no schema check)'
$ws.Range('C100').Value = 'motor_current_percent=''92.0'' and updateTime<''2021-05-21T14:22:43'' or id=''38'''
$ws.Range('F100').Value = 'HeatPumpKpiData'
$ws.Range('H100').Value = 0
$ws.Range('I100').Value = 0
$ws.Range('L100').Value = 0

$ws.Range('A101').Value = 'iems-connector-test-mysql-Float_varchar-String_varchar-Long_int-3'
$ws.Range('B101').Value = 'good request, data retrieved (no schema check)'
$ws.Range('C101').Value = 'motor_current_percent=''92.0'' or deviceName=''1#制冷机'' and id=''38'''
$ws.Range('F101').Value = 'HeatPumpKpiData'
$ws.Range('H101').Value = 0
$ws.Range('I101').Value = 0
$ws.Range('L101').Value = 0

$ws.Range('A102').Value = 'iems-connector-test-mysql-LocalDateTime_datetime-String_varchar-Long_int-3'
$ws.Range('B102').Value = 'good request, data retrieved (no schema check)'
$ws.Range('C102').Value = 'updateTime<''2021-05-21T14:22:43'' or deviceName=''1#制冷机'' and id=''38'''
$ws.Range('F102').Value = 'HeatPumpKpiData'
$ws.Range('H102').Value = 0
$ws.Range('I102').Value = 0
$ws.Range('L102').Value = 0

$ws.Range('A103').Value = 'iems-connector-test-mysql-Float_varchar-LocalDateTime_datetime-String_varchar-Long_int-3'
$ws.Range('B103').Value = 'good request, data retrieved (no schema check)'
$ws.Range('C103').Value = 'motor_current_percent=''92.0'' or updateTime<''2021-05-21T14:22:43'' and deviceName=''1#制冷机'' or id=''38'''
$ws.Range('F103').Value = 'HeatPumpKpiData'
$ws.Range('H103').Value = 0
$ws.Range('I103').Value = 0
$ws.Range('L103').Value = 0

# Copy cell formatting (style 3 = quote-prefixed text format) from an existing B-column cell onto the new B-column cells only
$ws.Range('B76').Copy() | Out-Null
$ws.Range('B77:B103').PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Widen column A to fit the new, longer test-case names
$ws.Columns("A").ColumnWidth = 53.59

# Move the selection / frozen-pane viewport down to the newly added rows
$ws.Range("A68").Select() | Out-Null
$ws.Range("A105").Select() | Out-Null